# Update cryptocurrency price/volume data, and fix two mis-ordered rows
# (RenderToken/PancakeSwap at rows 27-28, TheGraph/dogwifhat at rows 42/44)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.561.95"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "3.186.43"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'601.29"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'155.79"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.188.76"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "'0.548"
$ws.Range("E9").Value = "  +3.01%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'5.91"
$ws.Range("E11").Value = "  -4.37%  "
$ws.Range("D12").Value = "'0.510"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'0.0000264"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "'38.84"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "3.721.28"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "66.658.74"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "'7.38"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "3.193.66"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'513.99"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "'15.35"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'8.12"
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("D24").Value = "'14.88"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'84.80"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'9.23"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'2.40"
$ws.Range("E29").Value = "  +8.35%  "
$ws.Range("D30").Value = "'3.06"
$ws.Range("E30").Value = "  +7.38%  "
$ws.Range("D31").Value = "'7.00"
$ws.Range("E31").Value = "  +7.31%  "
$ws.Range("D32").Value = "'28.03"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "'6.52"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").Value = "'513.12"
$ws.Range("E36").Value = "  +7.76%  "
$ws.Range("D37").Value = "'54.89"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'0.0893"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "'0.0421"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'8.86"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'0.124"
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.303"
$ws.Range("E42").Value = "  +5.63%  "
$ws.Range("D43").Value = "0.0₃0681"
$ws.Range("E43").Value = "  +11.00%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "2.858.80"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("D47").Value = "'28.38"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  +5.11%  "
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "'2.64"
$ws.Range("E51").Value = "  +8.13%  "
